# Append new content to the end of the third paragraph (the one that ends
# with "...submit my homework)."), right before its paragraph mark.
$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(3)
$pr = $p.Range

# A Range that stops just before the paragraph mark (so InsertAfter appends
# new runs inside the paragraph, not after it).
$ip = $d.Range($pr.Start, $pr.End - 1)
$ip.Collapse(0)   # wdCollapseEnd

$lineBreak = [char]11   # wdLineBreak plain-text representation

# Two separate <w:br/> runs.
$ip.InsertAfter($lineBreak.ToString())
$ip.Collapse(0)
$ip.InsertAfter($lineBreak.ToString())
$ip.Collapse(0)

# "I" - own run (font-hint run in the source document; formatting left as-is).
$ip.InsertAfter("I")
$ip.Collapse(0)

# " believe I also did achieve something when I was doing my lab assignments. I successfully ran the 4"
$ip.InsertAfter(" believe I also did achieve something when I was doing my lab assignments. I successfully ran the 4")
$ip.Collapse(0)

# "th" superscript
$ip.InsertAfter("th")
$thRange = $d.Range($ip.Start, $ip.End)
$thRange.Font.Superscript = $true
$ip.Collapse(0)

# " question by installing several .jar class files and extensions into Apache "
$ip.InsertAfter(" question by installing several .jar class files and extensions into Apache ")
$ip.Collapse(0)

# "Netbeans"
$ip.InsertAfter("Netbeans")
$ip.Collapse(0)

# " application, and I was able to execute the code with an external output using Java. It's my first
# time seeing a visualized chart from a sole code program. I was amused by the power of coding, and it
# definitely motivates me to dive deeper into this world of programming, just like saying "Hello World"
# to a whole new universe."
$ip.InsertAfter(" application, and I was able to execute the code with an external output using Java. It" + [char]8217 + "s my first time seeing a visualized chart from a sole code program. I was amused by the power of coding, and it definitely motivates me to dive deeper into this world of programming, just like saying " + [char]8220 + "Hello World" + [char]8221 + " to a whole new universe.")
